$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.813.94"
$ws.Range("E2").Value = "  +2.75%  "

# Row 3
$ws.Range("D3").Value = "2.974.02"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'595.26"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$ws.Range("D6").Value = "'145.86"
$ws.Range("E6").Value = "  +0.52%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "2.970.95"
$ws.Range("E8").Value = "  +1.26%  "

# Row 9
$ws.Range("E9").Value = "  +0.42%  "

# Row 10
$ws.Range("D10").Value = "'7.24"
$ws.Range("E10").Value = "  +3.52%  "

# Row 11
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "  +3.87%  "

# Row 12
$ws.Range("E12").Value = "  +1.10%  "

# Row 13
$ws.Range("E13").Value = "  +5.82%  "

# Row 14
$ws.Range("D14").Value = "'33.28"
$ws.Range("E14").Value = "  -1.37%  "

# Row 15
$ws.Range("E15").Value = "  -0.39%  "

# Row 16
$ws.Range("D16").Value = "3.464.98"
$ws.Range("E16").Value = "  +1.37%  "

# Row 17
$ws.Range("D17").Value = "62.787.98"
$ws.Range("E17").Value = "  +2.84%  "

# Row 18
$ws.Range("D18").Value = "'6.72"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19
$ws.Range("D19").Value = "2.975.82"
$ws.Range("E19").Value = "  +1.55%  "

# Row 20
$ws.Range("D20").Value = "'443.39"
$ws.Range("E20").Value = "  +2.07%  "

# Row 21
$ws.Range("D21").Value = "'13.51"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
$ws.Range("E22").Value = "  -1.04%  "

# Row 23
$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'81.87"
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'11.27"
$ws.Range("E25").Value = "  +1.70%  "

# Row 26
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("E27").Value = "  -3.52%  "

# Row 28
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  +3.86%  "

# Row 30
$ws.Range("D30").Value = "'2.62"
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("E31").Value = "  -4.67%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0965"
$ws.Range("E32").Value = "  +10.90%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("D34").Value = "'26.54"
$ws.Range("E34").Value = "  -0.83%  "

# Row 35
$ws.Range("E35").Value = "  +0.16%  "

# Row 36
$ws.Range("D36").Value = "'0.996"
$ws.Range("E36").Value = "  -1.72%  "

# Row 37
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  +3.32%  "

# Row 39
$ws.Range("E39").Value = "  +3.10%  "

# Row 40
$ws.Range("D40").Value = "'49.55"
$ws.Range("E40").Value = "  -0.88%  "

# Row 41
$ws.Range("D41").Value = "'8.56"
$ws.Range("E41").Value = "  -0.46%  "

# Row 42
$ws.Range("E42").Value = "  -4.29%  "

# Row 43
$ws.Range("E43").Value = "  -0.61%  "

# Row 44
$ws.Range("D44").Value = "'40.23"
$ws.Range("E44").Value = "  -4.57%  "

# Row 45
$ws.Range("D45").Value = "2.749.05"
$ws.Range("E45").Value = "  +1.39%  "

# Row 46
$ws.Range("D46").Value = "'135.50"
$ws.Range("E46").Value = "  +1.62%  "

# Row 47
$ws.Range("E47").Value = "  -1.55%  "

# Row 48
$ws.Range("D48").Value = "'363.16"
$ws.Range("E48").Value = "  -2.67%  "

# Row 50
$ws.Range("D50").Value = "'23.13"
$ws.Range("E50").Value = "  -3.08%  "

# Row 51
$ws.Range("E51").Value = "  -0.10%  "
